$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Grupo A")
$ws.Range("H2").Value = 54.6
$ws.Range("I2").Value = 12.26
$ws.Range("G3").Value = 49
$ws.Range("H3").Value = 47.86
$ws.Range("I3").Value = 1.140000000000001
$ws.Range("G4").Value = 54.6
$ws.Range("I4").Value = -12.26
$ws.Range("G5").Value = 47.86
$ws.Range("H5").Value = 49
$ws.Range("I5").Value = -1.140000000000001

$ws = $wb.Worksheets.Item("Grupo B")
$ws.Range("G2").Value = 84.86
$ws.Range("H2").Value = 49.76
$ws.Range("I2").Value = 35.1
$ws.Range("G3").Value = 72.7
$ws.Range("H3").Value = 47.86
$ws.Range("I3").Value = 24.84
$ws.Range("G4").Value = 49.76
$ws.Range("H4").Value = 84.86
$ws.Range("I4").Value = -35.1
$ws.Range("G5").Value = 47.86
$ws.Range("H5").Value = 72.7
$ws.Range("I5").Value = -24.84

$ws = $wb.Worksheets.Item("Grupo C")
$ws.Range("G2").Value = 74.06
$ws.Range("H2").Value = 61.56
$ws.Range("I2").Value = 12.5
$ws.Range("G3").Value = 72.86
$ws.Range("H3").Value = 60.16
$ws.Range("I3").Value = 12.7
$ws.Range("G4").Value = 61.56
$ws.Range("H4").Value = 74.06
$ws.Range("I4").Value = -12.5
$ws.Range("G5").Value = 60.16
$ws.Range("H5").Value = 72.86
$ws.Range("I5").Value = -12.7

$ws = $wb.Worksheets.Item("Grupo D")
$ws.Range("G2").Value = 73.95999999999999
$ws.Range("H2").Value = 53.66
$ws.Range("I2").Value = 20.3
$ws.Range("G3").Value = 68.06
$ws.Range("H3").Value = 64.7
$ws.Range("I3").Value = 3.359999999999999
$ws.Range("G4").Value = 64.7
$ws.Range("H4").Value = 68.06
$ws.Range("I4").Value = -3.359999999999999
$ws.Range("G5").Value = 53.66
$ws.Range("H5").Value = 73.95999999999999
$ws.Range("I5").Value = -20.3

$ws = $wb.Worksheets.Item("Grupo E")
$ws.Range("G2").Value = 84.26000000000001
$ws.Range("H2").Value = 63.9
$ws.Range("I2").Value = 20.36000000000001
$ws.Range("G3").Value = 72.45
$ws.Range("H3").Value = 54.16
$ws.Range("I3").Value = 18.29000000000001
$ws.Range("G4").Value = 63.9
$ws.Range("H4").Value = 84.26000000000001
$ws.Range("I4").Value = -20.36000000000001
$ws.Range("G5").Value = 54.16
$ws.Range("H5").Value = 72.45
$ws.Range("I5").Value = -18.29000000000001

$ws = $wb.Worksheets.Item("Grupo F")
$ws.Range("G2").Value = 62.56
$ws.Range("H2").Value = 38.26
$ws.Range("B3").Value = 'DM Studio'
$ws.Range("G3").Value = 60.2
$ws.Range("H3").Value = 59.25
$ws.Range("I3").Value = 0.9500000000000028
$ws.Range("B4").Value = 'Rolo Compressor ZN'
$ws.Range("G4").Value = 59.25
$ws.Range("H4").Value = 60.2
$ws.Range("I4").Value = -0.9500000000000028
$ws.Range("G5").Value = 38.26
$ws.Range("H5").Value = 62.56

$ws = $wb.Worksheets.Item("Grupo G")
$ws.Range("G2").Value = 73.76000000000001
$ws.Range("H2").Value = 57.6
$ws.Range("I2").Value = 16.16
$ws.Range("G3").Value = 61.96
$ws.Range("H3").Value = 54.1
$ws.Range("I3").Value = 7.859999999999999
$ws.Range("B4").Value = 'Grêmio imortal 37'
$ws.Range("G4").Value = 57.6
$ws.Range("H4").Value = 73.76000000000001
$ws.Range("I4").Value = -16.16
$ws.Range("B5").Value = 'A Lenda Super Vascão f.c'
$ws.Range("G5").Value = 54.1
$ws.Range("H5").Value = 61.96
$ws.Range("I5").Value = -7.859999999999999

$ws = $wb.Worksheets.Item("Grupo H")
$ws.Range("G2").Value = 68.06
$ws.Range("H2").Value = 58.96
$ws.Range("G3").Value = 57.45
$ws.Range("H3").Value = 42.96
$ws.Range("I3").Value = 14.49
$ws.Range("G4").Value = 58.96
$ws.Range("H4").Value = 68.06
$ws.Range("G5").Value = 42.96
$ws.Range("H5").Value = 57.45
$ws.Range("I5").Value = -14.49
